$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Honza Vrátník")

# Update the text of A26 (shared string used for "3. iterace - modely balíků")
$ws.Range("A26").Value = "3. iterace - modely balíků, pár dalších scénářů"

# Update the hours value for that row (B26): 1 -> 1.5
$ws.Range("B26").Value = 1.5

# Move the active selection from B27 to A27
$ws.Range("A27").Select()
